# "Better plot for posters" - update the GLUE parameter bounds for row 14
# (B14/C14) and leave the sheet scrolled/selected near that row, matching
# what the author would have seen right after editing those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the parameter range used to drive the plot.
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 8

# Leave the selection on the last-edited cell, as captured in the saved view.
$ws.Range("D14").Select()
